# Weekly data refresh: insert this week's new price record at the top of the
# data block (row 116), pushing all existing rows down by one. The oldest
# record that falls off the bottom of the 58-row window becomes the new
# last row (174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 116; rows 116:173 shift down to 117:174.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(116, 1).Value = 4
$ws.Cells.Item(116, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(116, 3).Value = "Los Lagos"
$ws.Cells.Item(116, 4).Value = 44466
$ws.Cells.Item(116, 5).Value = 10
$ws.Cells.Item(116, 6).Value = 100112045
$ws.Cells.Item(116, 7).Value = "Zapallo"
$ws.Cells.Item(116, 8).Value = "Paine"
$ws.Cells.Item(116, 9).Value = "1a (guarda)"
$ws.Cells.Item(116, 10).Value = 500
$ws.Cells.Item(116, 11).Value = 600
$ws.Cells.Item(116, 12).Value = 600
$ws.Cells.Item(116, 13).Value = 600
$ws.Cells.Item(116, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 600
$ws.Cells.Item(116, 17).Value = 1
$ws.Cells.Item(116, 18).Value = "Hortaliza"
